$d = $word.ActiveDocument
$d.Content.Find.Execute("future  career", $false, $false, $false, $false, $false, $true, 1, $false, "future career", 2)
